# Generate Report for Handoff
# Update status text and timestamps across the Overview, zh-cn and de-de
# sheets, then re-autofit the now-shorter "Status" columns.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# New status text (was: "Handed back: in sync with en-US")
$newStatus = "Ready for handoff"

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-13 21:16:59"

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-13 21:16:51"

$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-13 21:16:59"

# Re-fit the Status columns now that the text is shorter (~17.2 chars wide).
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
